$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin (B), Link (C), Price (D) and Volume(1h) (E) columns for rows 2-51
# to reflect the refreshed cryptos list.

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "27.031.93"
$ws.Range("E2").Value = "  -3.00%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.719.07"
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Formula = "=""1.004"""
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Formula = "=""316.18"""
$ws.Range("E5").Value = "  -3.58%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Formula = "=""1.003"""
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Formula = "=""0.4616"""
$ws.Range("E7").Value = "  +3.01%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Formula = "=""0.3444"""
$ws.Range("E8").Value = "  -3.57%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Formula = "=""42.97"""
$ws.Range("E9").Value = "  +1.94%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Formula = "=""0.07309"""
$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Formula = "=""1.051"""
$ws.Range("E11").Value = "  -3.96%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").Formula = "=""1.003"""
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Formula = "=""19.87"""
$ws.Range("E13").Value = "  -4.80%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Formula = "=""5.872"""
$ws.Range("E14").Value = "  -2.93%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.710.26"
$ws.Range("E15").Value = "  -3.99%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Formula = "=""6.916"""
$ws.Range("E16").Value = "  -4.15%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Formula = "=""89.41"""
$ws.Range("E17").Value = "  -3.85%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Formula = "=""0.00001046"""
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Formula = "=""0.06332"""
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Formula = "=""1.004"""
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Formula = "=""16.53"""
$ws.Range("E21").Value = "  -3.83%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Formula = "=""5.636"""
$ws.Range("E22").Value = "  -3.19%  "

$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "27.073.68"
$ws.Range("E23").Value = "  -2.90%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Formula = "=""10.85"""
$ws.Range("E24").Value = "  -4.20%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Formula = "=""2.153"""
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Formula = "=""157.23"""
$ws.Range("E26").Value = "  -3.28%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Formula = "=""19.45"""
$ws.Range("E27").Value = "  -3.90%  "

$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "1.912.27"
$ws.Range("E28").Value = "  -3.51%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Formula = "=""2.137"""
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Formula = "=""119.24"""
$ws.Range("E30").Value = "  -4.91%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Formula = "=""1.027"""
$ws.Range("E31").Value = "  -6.78%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Formula = "=""0.09108"""
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Formula = "=""3.592"""
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Formula = "=""5.331"""
$ws.Range("E34").Value = "  -4.59%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Formula = "=""0.02203"""
$ws.Range("E35").Value = "  -3.96%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Formula = "=""11.15"""
$ws.Range("E36").Value = "  -5.75%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Formula = "=""0.05852"""
$ws.Range("E37").Value = "  -3.88%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Formula = "=""4.769"""
$ws.Range("E38").Value = "  -3.76%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Formula = "=""0.1995"""
$ws.Range("E39").Value = "  -4.93%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Formula = "=""1.404"""
$ws.Range("E40").Value = "  +0.74%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Formula = "=""0.5968"""
$ws.Range("E41").Value = "  -5.72%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Formula = "=""1.126"""
$ws.Range("E42").Value = "  -4.68%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Formula = "=""7.526"""
$ws.Range("E43").Value = "  -5.04%  "

$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").Formula = "=""3.627"""
$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "=""12.61"""
$ws.Range("E45").Value = "  -4.56%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Formula = "=""0.5636"""
$ws.Range("E46").Value = "  -4.03%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Formula = "=""119.68"""
$ws.Range("E47").Value = "  -2.38%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Formula = "=""1.866"""
$ws.Range("E48").Value = "  -4.63%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Formula = "=""0.06654"""
$ws.Range("E49").Value = "  -3.68%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Formula = "=""1.088"""
$ws.Range("E50").Value = "  -4.46%  "

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Formula = "=""1.003"""
$ws.Range("E51").Value = "  +0.06%  "

# Convert the temporary text-formulas in column D into static values
# (Copy + PasteSpecial values-only keeps them typed as text, matching
# the original inlineStr cells, without touching cell styles).
$ws.Range("D4:D14").Copy()
$ws.Range("D4:D14").PasteSpecial(-4163)
$ws.Range("D16:D22").Copy()
$ws.Range("D16:D22").PasteSpecial(-4163)
$ws.Range("D24:D27").Copy()
$ws.Range("D24:D27").PasteSpecial(-4163)
$ws.Range("D29:D51").Copy()
$ws.Range("D29:D51").PasteSpecial(-4163)
$excel.CutCopyMode = 0
